$wb = $excel.ActiveWorkbook

# Select the CRtPaL-profits sheet and update values B2:B24 from 1 to 2
$wsProfits = $wb.Worksheets.Item("CRtPaL-profits")
$wsProfits.Select()

$wsProfits.Range("B2:B24").Value = 2
$wsProfits.Range("B2:B24").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Select()
